$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 994.4
$ws.Range("I129").Value = 560.6667
$ws.Range("J129").Value = 1180.2858
$ws.Range("K129").Value = 1682.0001
$ws.Range("L129").Value = 3540.8574
$ws.Range("M129").Value = 3317.9999
$ws.Range("N129").Value = -13540.8574

$ws.Range("H135").Value = 40000756
$ws.Range("I135").Value = 788.0417
$ws.Range("J135").Value = 1000000000
$ws.Range("K135").Value = 7092.3753
$ws.Range("L135").Value = 9000000000
$ws.Range("M135").Value = -4557.3753
$ws.Range("N135").Value = -9000005070

$ws.Range("H137").Value = 2223947
$ws.Range("I137").Value = 1446.25
$ws.Range("J137").Value = 6175059.5
$ws.Range("K137").Value = 4338.75
$ws.Range("L137").Value = 18525178.5
$ws.Range("M137").Value = -1788.75
$ws.Range("N137").Value = -18530278.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4173.24
$ws.Range("I32").Value = 4054.8867
$ws.Range("J32").Value = 8000
$ws.Range("K32").Value = 4054.8867
$ws.Range("L32").Value = 8000
$ws.Range("M32").Value = -3767.8867
$ws.Range("N32").Value = -8574

$ws.Range("H43").Value = 6906
$ws.Range("J43").Value = 6906
$ws.Range("L43").Value = 6906
$ws.Range("N43").Value = -7532

$ws.Range("H92").Value = 27777.5
$ws.Range("J92").Value = 27777.5
$ws.Range("L92").Value = 27777.5
$ws.Range("N92").Value = -32769.5

$ws.Range("H109").Value = 30000
$ws.Range("J109").Value = 30000
$ws.Range("L109").Value = 30000
$ws.Range("N109").Value = -32774

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H86").Value = 1619.2693
$ws.Range("I86").Value = 1536.409
$ws.Range("J86").Value = 2075
$ws.Range("K86").Value = 1536.409
$ws.Range("L86").Value = 2075
$ws.Range("M86").Value = -413.4090000000001
$ws.Range("N86").Value = -4321

$ws.Range("H89").Value = 1619.2693
$ws.Range("I89").Value = 1536.409
$ws.Range("J89").Value = 2075
$ws.Range("K89").Value = 7682.045
$ws.Range("L89").Value = 10375
$ws.Range("M89").Value = -2066.045
$ws.Range("N89").Value = -21607

$ws.Range("H99").Value = 1593.4117
$ws.Range("I99").Value = 1336.1818
$ws.Range("J99").Value = 2065
$ws.Range("K99").Value = 1336.1818
$ws.Range("L99").Value = 2065
$ws.Range("M99").Value = 161.8181999999999
$ws.Range("N99").Value = -5061

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 23840.354
$ws.Range("I6").Value = 26619.066
$ws.Range("J6").Value = 3000
$ws.Range("K6").Value = 26619.066
$ws.Range("L6").Value = 3000
$ws.Range("M6").Value = -26506.066
$ws.Range("N6").Value = -3226

$ws.Range("H31").Value = 2181.6667
$ws.Range("I31").Value = 1412
$ws.Range("J31").Value = 2896.3572
$ws.Range("K31").Value = 1412
$ws.Range("L31").Value = 2896.3572
$ws.Range("M31").Value = -1117
$ws.Range("N31").Value = -3486.3572

$ws.Range("H33").Value = 18250
$ws.Range("I33").Value = 5333.3335
$ws.Range("J33").Value = 26000
$ws.Range("K33").Value = 5333.3335
$ws.Range("L33").Value = 26000
$ws.Range("M33").Value = -4954.3335
$ws.Range("N33").Value = -26758

$ws.Range("H34").Value = 2181.6667
$ws.Range("I34").Value = 1412
$ws.Range("J34").Value = 2896.3572
$ws.Range("K34").Value = 1412
$ws.Range("L34").Value = 2896.3572
$ws.Range("M34").Value = -1210
$ws.Range("N34").Value = -3300.3572

$ws.Range("H62").Value = 2411.111
$ws.Range("I62").Value = 2416.6667
$ws.Range("J62").Value = 2400
$ws.Range("K62").Value = 2416.6667
$ws.Range("L62").Value = 2400
$ws.Range("M62").Value = -1792.6667
$ws.Range("N62").Value = -3648

$ws.Range("H65").Value = 2411.111
$ws.Range("I65").Value = 2416.6667
$ws.Range("J65").Value = 2400
$ws.Range("K65").Value = 12083.3335
$ws.Range("L65").Value = 12000
$ws.Range("M65").Value = -8963.333500000001
$ws.Range("N65").Value = -18240

$ws.Range("H99").Value = 3318.2
$ws.Range("I99").Value = 2772.75
$ws.Range("K99").Value = 2772.75
$ws.Range("M99").Value = -1274.75

$ws.Range("H126").Value = 3318.2
$ws.Range("I126").Value = 2772.75
$ws.Range("K126").Value = 8318.25
$ws.Range("M126").Value = -5848.25

$ws.Range("H134").Value = 1743.9615
$ws.Range("I134").Value = 1423.6086
$ws.Range("J134").Value = 4200
$ws.Range("K134").Value = 4270.825800000001
$ws.Range("L134").Value = 12600
$ws.Range("M134").Value = -1735.825800000001
$ws.Range("N134").Value = -17670

$ws.Range("H141").Value = 49208.273
$ws.Range("J141").Value = 49208.273
$ws.Range("L141").Value = 49208.273
$ws.Range("N141").Value = -59568.273

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 119.22222
$ws.Range("I23").Value = 86
$ws.Range("J23").Value = 145.8
$ws.Range("K23").Value = 258
$ws.Range("L23").Value = 437.4
$ws.Range("M23").Value = -23
$ws.Range("N23").Value = -907.4000000000001

$ws.Range("H122").Value = 21272.623
$ws.Range("I122").Value = 30123.676
$ws.Range("J122").Value = 804.5625
$ws.Range("K122").Value = 271113.084
$ws.Range("L122").Value = 7241.0625
$ws.Range("M122").Value = -268663.084
$ws.Range("N122").Value = -12141.0625

$ws.Range("H130").Value = 1730.3334
$ws.Range("I130").Value = 1337.25
$ws.Range("J130").Value = 2516.5
$ws.Range("K130").Value = 4011.75
$ws.Range("L130").Value = 7549.5
$ws.Range("M130").Value = 1008.25
$ws.Range("N130").Value = -17589.5

$ws.Range("H131").Value = 1772.0526
$ws.Range("I131").Value = 2388.2
$ws.Range("J131").Value = 1737.8223
$ws.Range("K131").Value = 7164.599999999999
$ws.Range("L131").Value = 5213.4669
$ws.Range("M131").Value = -2124.599999999999
$ws.Range("N131").Value = -15293.4669

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 14000
$ws.Range("J26").Value = 14000
$ws.Range("L26").Value = 14000
$ws.Range("N26").Value = -14560

$ws.Range("H50").Value = 14000
$ws.Range("J50").Value = 14000
$ws.Range("L50").Value = 14000
$ws.Range("N50").Value = -14996

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 35000
$ws.Range("J39").Value = 35000
$ws.Range("L39").Value = 35000
$ws.Range("N39").Value = -35920

$ws.Range("H44").Value = 10000
$ws.Range("J44").Value = 10000
$ws.Range("L44").Value = 10000
$ws.Range("N44").Value = -10912

$ws.Range("H56").Value = 4051
$ws.Range("I56").Value = 4051
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 4051
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -3360
$ws.Range("N56").ClearContents()

$ws.Range("H132").Value = 4905.778
$ws.Range("I132").Value = 5319.6665
$ws.Range("J132").Value = 4078
$ws.Range("K132").Value = 15958.9995
$ws.Range("L132").Value = 12234
$ws.Range("M132").Value = -13428.9995
$ws.Range("N132").Value = -17294

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 19600
$ws.Range("I58").Value = 11500
$ws.Range("J58").Value = 25000
$ws.Range("K58").Value = 11500
$ws.Range("L58").Value = 25000
$ws.Range("M58").Value = -11192
$ws.Range("N58").Value = -25616

$ws.Range("H61").Value = 18422.8
$ws.Range("I61").Value = 14000
$ws.Range("J61").Value = 25057
$ws.Range("K61").Value = 14000
$ws.Range("L61").Value = 25057
$ws.Range("M61").Value = -13708
